$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA2").Value = 2.3246209
$ws.Range("AH2").Value = 59
$ws.Range("AI2").Value = 112
$ws.Range("DF2").Value = 73
$ws.Range("DG2").Value = 174
$ws.Range("AZ3").Value = 33.333333333333
$ws.Range("BB3").Value = 66.666666666667
$ws.Range("BD3").Value = 63.636363636364
$ws.Range("CS3").Value = 6
$ws.Range("CU3").Value = 4
$ws.Range("K4").Value = 4.7902
$ws.Range("AA4").Value = 2.08300496
$ws.Range("BB4").Value = 47.619047619048
$ws.Range("BF4").Value = 51.546391752577
$ws.Range("CU4").Value = 99
$ws.Range("CV4").Value = 47
$ws.Range("AA5").Value = 0.0753351
$ws.Range("AE5").Value = 134
$ws.Range("AF5").Value = 79.761904761905
$ws.Range("AH5").Value = 61
$ws.Range("BG5").Value = 46
$ws.Range("BT5").Value = 34
$ws.Range("K7").Value = 1.6361
$ws.Range("AA7").Value = 5.25996753
$ws.Range("CC7").Value = 85
$ws.Range("AA8").Value = 2.01561185
$ws.Range("AB8").Value = 691
$ws.Range("AS8").Value = 24
$ws.Range("AY8").Value = 7
$ws.Range("AZ8").Value = 46.666666666667
$ws.Range("BA8").Value = 75
$ws.Range("BB8").Value = 44.910179640719
$ws.Range("BC8").Value = 56
$ws.Range("BD8").Value = 44.094488188976
$ws.Range("CS8").Value = 15
$ws.Range("CU8").Value = 92
$ws.Range("K9").Value = 1.7392
$ws.Range("AA10").Value = 1.40573129
$ws.Range("AB10").Value = 883
$ws.Range("AE10").Value = 599
$ws.Range("AF10").Value = 87.44525547445301
$ws.Range("AG10").Value = 685
$ws.Range("AH10").Value = 270
$ws.Range("DF10").Value = 289
$ws.Range("AA11").Value = 0.64804127
$ws.Range("K12").Value = 0.5459000000000001
$ws.Range("AA12").Value = 0.37221976
$ws.Range("AS12").Value = 30
$ws.Range("BB12").Value = 49.673202614379
$ws.Range("BD12").Value = 48.031496062992
$ws.Range("CU12").Value = 77
$ws.Range("AA13").Value = 1.22050169
$ws.Range("AJ13").Value = 34
$ws.Range("AA15").Value = 1.4622784
$ws.Range("AB15").Value = 1211
$ws.Range("AE15").Value = 652
$ws.Range("AF15").Value = 82.32323232323201
$ws.Range("AG15").Value = 792
$ws.Range("AI15").Value = 262
$ws.Range("AP15").Value = 46
$ws.Range("BA15").Value = 79
$ws.Range("BB15").Value = 55.244755244755
$ws.Range("BC15").Value = 64
$ws.Range("BD15").Value = 58.715596330275
$ws.Range("DB15").Value = 76.086956521739
$ws.Range("DG15").Value = 351
$ws.Range("K16").Value = 1.9401
$ws.Range("AA16").Value = 2.40439832
$ws.Range("AS16").Value = 23
$ws.Range("BB16").Value = 51.834862385321
$ws.Range("BD16").Value = 53.030303030303
$ws.Range("CU16").Value = 105
$ws.Range("DA16").Value = 42
$ws.Range("DB16").Value = 61.764705882353
$ws.Range("K17").Value = 2.2415
$ws.Range("AA17").Value = 1.12583059
$ws.Range("AA18").Value = 0.12204902
$ws.Range("AB18").Value = 554
$ws.Range("BG18").Value = 57
$ws.Range("AA23").Value = 0.12427912
$ws.Range("DE23").Value = 18
$ws.Range("DL23").Value = -2.6253
